$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as the most recent entry (row 109),
# pushing all subsequent rows (old 109..165) down by one (new 110..166).
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new record's data.
$ws.Cells.Item(109, 1).Value = 4
$ws.Cells.Item(109, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(109, 3).Value = "Los Lagos"
$ws.Cells.Item(109, 4).Value = 44960
$ws.Cells.Item(109, 5).Value = 10
$ws.Cells.Item(109, 6).Value = 100112052
$ws.Cells.Item(109, 7).Value = "Albahaca"
$ws.Cells.Item(109, 8).Value = "Sin especificar"
$ws.Cells.Item(109, 9).Value = "Primera"
$ws.Cells.Item(109, 10).Value = 160
$ws.Cells.Item(109, 11).Value = 5000
$ws.Cells.Item(109, 12).Value = 6000
$ws.Cells.Item(109, 13).Value = 5500
$ws.Cells.Item(109, 14).Value = "$/docena de matas"
$ws.Cells.Item(109, 15).Value = "Región Metropolitana"
$ws.Cells.Item(109, 16).Value = 917
$ws.Cells.Item(109, 17).Value = 6
$ws.Cells.Item(109, 18).Value = "Hortaliza"
